$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 50 with style copied from row 49 (column A uses the bordered/bold style)
$ws.Range("A49").Copy() | Out-Null
$ws.Range("A50").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(50, 1).Value = 48

# Update B2:C50 values (recomputed via priority-queue based algorithm)
$ws.Cells.Item(2, 2).Value = 2.465091977227747
$ws.Cells.Item(2, 3).Value = 3.692498141701704
$ws.Cells.Item(3, 2).Value = 3.98674660142422
$ws.Cells.Item(3, 3).Value = 8.257120541548455
$ws.Cells.Item(4, 2).Value = 6.320821553804871
$ws.Cells.Item(4, 3).Value = 12.58096189793132
$ws.Cells.Item(5, 2).Value = 7.295785279507488
$ws.Cells.Item(5, 3).Value = 16.6678810491552
$ws.Cells.Item(6, 2).Value = 9.296406709044936
$ws.Cells.Item(6, 3).Value = 20.48535161316388
$ws.Cells.Item(7, 2).Value = 11.40254855409558
$ws.Cells.Item(7, 3).Value = 24.53397119679713
$ws.Cells.Item(8, 2).Value = 13.62433808955116
$ws.Cells.Item(8, 3).Value = 28.74625142898509
$ws.Cells.Item(9, 2).Value = 15.81063297057566
$ws.Cells.Item(9, 3).Value = 32.79532050261079
$ws.Cells.Item(10, 2).Value = 16.65153234147942
$ws.Cells.Item(10, 3).Value = 36.60154096709385
$ws.Cells.Item(11, 2).Value = 17.42086202135393
$ws.Cells.Item(11, 3).Value = 41.01487656732795
$ws.Cells.Item(12, 2).Value = 20.29479243278586
$ws.Cells.Item(12, 3).Value = 45.47795098574854
$ws.Cells.Item(13, 2).Value = 22.43181099330758
$ws.Cells.Item(13, 3).Value = 49.34874378640399
$ws.Cells.Item(14, 2).Value = 23.37899185267504
$ws.Cells.Item(14, 3).Value = 53.2933518821994
$ws.Cells.Item(15, 2).Value = 26.13930710349315
$ws.Cells.Item(15, 3).Value = 57.62198531478955
$ws.Cells.Item(16, 2).Value = 29.59143846964884
$ws.Cells.Item(16, 3).Value = 61.97989980019731
$ws.Cells.Item(17, 2).Value = 31.83510712643087
$ws.Cells.Item(17, 3).Value = 65.99720092520941
$ws.Cells.Item(18, 2).Value = 33.41030340071416
$ws.Cells.Item(18, 3).Value = 70.25700551979301
$ws.Cells.Item(19, 2).Value = 36.24794953173991
$ws.Cells.Item(19, 3).Value = 74.35407449972459
$ws.Cells.Item(20, 2).Value = 38.22688281316102
$ws.Cells.Item(20, 3).Value = 78.60941169504417
$ws.Cells.Item(21, 2).Value = 40.36573812032713
$ws.Cells.Item(21, 3).Value = 82.70924243515746
$ws.Cells.Item(22, 2).Value = 41.53970161309928
$ws.Cells.Item(22, 3).Value = 86.9816743317766
$ws.Cells.Item(23, 2).Value = 43.21082538763684
$ws.Cells.Item(23, 3).Value = 90.79920462701058
$ws.Cells.Item(24, 2).Value = 45.22656396597057
$ws.Cells.Item(24, 3).Value = 95.13090259869857
$ws.Cells.Item(25, 2).Value = 46.16781786297612
$ws.Cells.Item(25, 3).Value = 98.91194471471545
$ws.Cells.Item(26, 2).Value = 49.02109054520594
$ws.Cells.Item(26, 3).Value = 103.1482166925268
$ws.Cells.Item(27, 2).Value = 49.8586144872416
$ws.Cells.Item(27, 3).Value = 107.6629477656295
$ws.Cells.Item(28, 2).Value = 51.99609072055997
$ws.Cells.Item(28, 3).Value = 111.8702576416018
$ws.Cells.Item(29, 2).Value = 55.92027498799172
$ws.Cells.Item(29, 3).Value = 117.6219403525507
$ws.Cells.Item(30, 2).Value = 59.70210892696766
$ws.Cells.Item(30, 3).Value = 121.3400468714738
$ws.Cells.Item(31, 2).Value = 61.20959529930025
$ws.Cells.Item(31, 3).Value = 125.5519716545458
$ws.Cells.Item(32, 2).Value = 62.99657133408933
$ws.Cells.Item(32, 3).Value = 129.6518667243479
$ws.Cells.Item(33, 2).Value = 64.81888372075333
$ws.Cells.Item(33, 3).Value = 133.5378243134733
$ws.Cells.Item(34, 2).Value = 66.11370085601297
$ws.Cells.Item(34, 3).Value = 138.0309667284715
$ws.Cells.Item(35, 2).Value = 70.7375424185368
$ws.Cells.Item(35, 3).Value = 141.8989012554992
$ws.Cells.Item(36, 2).Value = 72.42480208725435
$ws.Cells.Item(36, 3).Value = 146.1518802730397
$ws.Cells.Item(37, 2).Value = 74.29362901006392
$ws.Cells.Item(37, 3).Value = 149.7598442004259
$ws.Cells.Item(38, 2).Value = 75.9025392317176
$ws.Cells.Item(38, 3).Value = 154.32911828045
$ws.Cells.Item(39, 2).Value = 77.34403965928399
$ws.Cells.Item(39, 3).Value = 158.7521628496195
$ws.Cells.Item(40, 2).Value = 77.80151278097104
$ws.Cells.Item(40, 3).Value = 162.6067022047367
$ws.Cells.Item(41, 2).Value = 79.25494671760097
$ws.Cells.Item(41, 3).Value = 167.728127799157
$ws.Cells.Item(42, 2).Value = 81.2554269542062
$ws.Cells.Item(42, 3).Value = 171.794438868152
$ws.Cells.Item(43, 2).Value = 82.78332659661285
$ws.Cells.Item(43, 3).Value = 176.0879111341021
$ws.Cells.Item(44, 2).Value = 84.30725099895773
$ws.Cells.Item(44, 3).Value = 180.5252297022838
$ws.Cells.Item(45, 2).Value = 86.73446636015788
$ws.Cells.Item(45, 3).Value = 185.1946300189942
$ws.Cells.Item(46, 2).Value = 88.25931581463048
$ws.Cells.Item(46, 3).Value = 189.6237043016379
$ws.Cells.Item(47, 2).Value = 89.56888375214771
$ws.Cells.Item(47, 3).Value = 193.8393638159411
$ws.Cells.Item(48, 2).Value = 93.19519895409421
$ws.Cells.Item(48, 3).Value = 197.7229260052057
$ws.Cells.Item(49, 2).Value = 95.61472246210842
$ws.Cells.Item(49, 3).Value = 202.9526369504172
$ws.Cells.Item(50, 2).Value = 97.40484373025957
$ws.Cells.Item(50, 3).Value = 207.1380315624927

Write-Output "done"
